$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New program: TUBerlin_ME
$ws.Cells.Item(9, 1).Value = "TUBerlin_ME"

# Existing programs (rows 2-8) switch their "Choose" value from Yes to No
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = "No"
}

$ws.Cells.Item(9, 2).Value = "Yes"

# New program: RWTH Aachen_ME
$ws.Cells.Item(10, 1).Value = "RWTH Aachen_ME"
$ws.Cells.Item(10, 2).Value = "Yes"

# New program: TUBraunschweig_ME
$ws.Cells.Item(11, 1).Value = "TUBraunschweig_ME"
$ws.Cells.Item(11, 2).Value = "Yes"

# Extend the Yes/No list validation to cover the newly added rows
$ws.Range("B1:B8").Validation.Delete()
$ws.Range("B1:B11").Validation.Add(3, 1, 3, '"Yes,No"')

$ws.Range("A12").Select()
